# Update gh-pages to output generated at 456a3b4
# Updates the "想去人数" (F) and "最低票价" (G) columns on the
# "展览" and "全部类型" worksheets, which hold duplicate data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Row => (F new value or $null, G new value or $null)
$updates = @{
    2  = @{ F = 110;  G = 40 }
    3  = @{          G = 70 }
    4  = @{          G = 60 }
    6  = @{          G = 60 }
    7  = @{ F = 467;  G = 60 }
    9  = @{ F = 194 }
    13 = @{ F = 290 }
    16 = @{ F = 1791 }
    18 = @{ F = 111 }
    23 = @{ F = 4330 }
    25 = @{ F = 306 }
    26 = @{ F = 1149 }
    27 = @{ F = 503 }
    29 = @{ F = 688 }
    31 = @{ F = 346 }
    33 = @{ F = 179 }
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $vals = $updates[$row]
        if ($vals.ContainsKey("F")) {
            $ws.Range("F$row").Value = $vals["F"]
        }
        if ($vals.ContainsKey("G")) {
            $ws.Range("G$row").Value = $vals["G"]
        }
    }
}
